# Weekly update: shift existing "Acelga" price records down by one
# week-pair (2 rows) to make room for a new week's data at the top of
# this date range (rows 214-215), and append the previous last record
# (old rows 254-255) as new rows 256-257.
#
# Columns that vary per record: D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg).
# All other columns (A,B,C,E,F,G,H,I,N,O,Q,R) repeat with a fixed
# Primera/Segunda pattern and do not need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colsToShift = @(4, 10, 11, 12, 13, 16)  # D, J, K, L, M, P

# 1) Append two brand-new rows (256, 257) that duplicate what used to be
#    the last pair of rows (254, 255), before any values are overwritten.
for ($col = 1; $col -le 18; $col++) {
    $v = $ws.Cells.Item(254, $col).Value2()
    $ws.Cells.Item(256, $col).Value = $v
    $v2 = $ws.Cells.Item(255, $col).Value2()
    $ws.Cells.Item(257, $col).Value = $v2
}
$ws.Cells.Item(256, 4).NumberFormat = $ws.Cells.Item(254, 4).NumberFormat
$ws.Cells.Item(257, 4).NumberFormat = $ws.Cells.Item(255, 4).NumberFormat

# 2) Shift rows 216..255 down from rows 214..253 (i.e. new row n = old row
#    n-2), working from the bottom up so sources are read before they are
#    overwritten.
for ($n = 255; $n -ge 216; $n--) {
    $src = $n - 2
    foreach ($col in $colsToShift) {
        $v = $ws.Cells.Item($src, $col).Value2()
        $ws.Cells.Item($n, $col).Value = $v
    }
}

# 3) Rows 214/215 become the newest week's data: new date (44505) and new
#    volume, keeping the same price range as before.
$ws.Cells.Item(214, 4).Value = 44505
$ws.Cells.Item(214, 10).Value = 2800

$ws.Cells.Item(215, 4).Value = 44505
$ws.Cells.Item(215, 10).Value = 1200
